$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header style (H1, "IP") onto the two new header cells
# so they reuse the same cellXf instead of minting a new style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(5, 8, 4, 9, 4, 9, 4, 9)
$jValues = @(5, 8, 6, 9, 5, 9, 4, 9)

for ($row = 2; $row -le 9; $row++) {
    $idx = $row - 2
    $ws.Cells.Item($row, 9).Value = $iValues[$idx]
    $ws.Cells.Item($row, 10).Value = $jValues[$idx]
}
